# Geschwindigkeitsberechnung - restructure Tabelle1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Wipe old layout completely (content + formatting) -----------------
$ws.Cells.Clear()

# --- Column widths (character units, converted so the saved XML width
#     matches the target: engine adds 5/6 to whatever ColumnWidth we set) --
$ws.Columns.Item(1).ColumnWidth = 2.1666666666666665     # A  -> 3
$ws.Columns.Item(2).ColumnWidth = 12.498697916666666      # B  -> 13.33203125
$ws.Columns.Item(3).ColumnWidth = 16.330729166666668      # C  -> 17.1640625
$ws.Columns.Item(4).ColumnWidth = 7.666666666666667        # D  -> 8.5
$ws.Columns.Item(5).ColumnWidth = 7.666666666666667        # E  -> 8.5
$ws.Columns.Item(6).ColumnWidth = 7.666666666666667        # F  -> 8.5
$ws.Columns.Item(7).ColumnWidth = 7.666666666666667        # G  -> 8.5
$ws.Columns.Item(8).ColumnWidth = 7.666666666666667        # H  -> 8.5
$ws.Columns.Item(9).ColumnWidth = 12.498697916666666       # I  -> 13.33203125
$ws.Columns.Item(10).ColumnWidth = 13.998697916666666      # J  -> 14.83203125
$ws.Columns.Item(11).ColumnWidth = 13.998697916666666      # K  -> 14.83203125

# --- Row heights -----------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 35
$ws.Rows.Item(3).RowHeight = 35
$ws.Rows.Item(4).RowHeight = 35
$ws.Rows.Item(5).RowHeight = 35

# --- Header row 2 ------------------------------------------------------
$ws.Range("B2").Value = "Distanz (m)"
$ws.Range("C2").Value = "Geschwindigkeit im Programm"
$ws.Range("D2").Value = "Zeitmessung (s)"
$ws.Range("J2").Value = "Reale Geschwindigkeit (km/h)"
$ws.Range("K2").Value = "Skalierung"

# --- Sub-header row 3 ----------------------------------------------------
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = "Durchschnitt"

# --- Data row 4 ------------------------------------------------------------
$ws.Range("B4").Value = 500
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 24.7
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 24.9
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 24.9
$ws.Range("I4").Formula = "=AVERAGE(D4:H4)"
$ws.Range("J4").Formula = "=B4/I4*3.6"
$ws.Range("K4").Formula = "=J4/C4"

# --- Data row 5 ------------------------------------------------------------
$ws.Range("B5").Value = 500
$ws.Range("C5").Value = 175
$ws.Range("D5").Value = 14.2
$ws.Range("E5").Value = 14.2
$ws.Range("F5").Value = 14.2
$ws.Range("G5").Value = 14.2
$ws.Range("H5").Value = 14.2
$ws.Range("I5").Formula = "=AVERAGE(D5:H5)"
$ws.Range("J5").Formula = "=B5/I5*3.6"
$ws.Range("K5").Formula = "=J5/C5"

# --- Styling (MUST happen before Merge(), so every covered cell - not just
#     the merge anchor - gets the style baked into the saved XML) ---------

function Style-DarkHeader($r) {
  $r.Font.Bold = $true
  $r.Interior.Pattern = 1
  $r.Interior.ThemeColor = 2
  $r.Interior.TintAndShade = -0.34998626667073579
  $r.Borders.LineStyle = 1
  $r.HorizontalAlignment = -4108
  $r.VerticalAlignment = -4108
  $r.WrapText = $true
}

function Style-LightHeader($r) {
  $r.Font.Bold = $true
  $r.Interior.Pattern = 1
  $r.Interior.ThemeColor = 2
  $r.Interior.TintAndShade = -0.14999847407452621
  $r.Borders.LineStyle = 1
  $r.HorizontalAlignment = -4108
  $r.VerticalAlignment = -4108
  $r.WrapText = $true
}

Style-DarkHeader $ws.Range("B2:B3")
Style-DarkHeader $ws.Range("C2:C3")
Style-DarkHeader $ws.Range("D2:I2")
Style-DarkHeader $ws.Range("J2:J3")
Style-DarkHeader $ws.Range("K2:K3")

Style-LightHeader $ws.Range("D3:I3")

# Data cells B4:J5 - bordered, centered, no fill, not bold
$dataCells = $ws.Range("B4:J5")
$dataCells.Borders.LineStyle = 1
$dataCells.HorizontalAlignment = -4108
$dataCells.VerticalAlignment = -4108

# Scaling column K4:K5 - bordered, centered, bold
$scaleCells = $ws.Range("K4:K5")
$scaleCells.Font.Bold = $true
$scaleCells.Borders.LineStyle = 1
$scaleCells.HorizontalAlignment = -4108
$scaleCells.VerticalAlignment = -4108

# --- Merges (after styling) ---------------------------------------------
$ws.Range("B2:B3").Merge()
$ws.Range("C2:C3").Merge()
$ws.Range("D2:I2").Merge()
$ws.Range("J2:J3").Merge()
$ws.Range("K2:K3").Merge()

# --- Sheet-level view / page setup --------------------------------------
$ws.Range("F16").Select()
$ws.PageSetup.PaperSize = 9

# --- Ignore the "formula omits adjacent cells" warning on I4:I5 --------
try {
  $ws.Range("I4:I5").Errors.Item(1).Ignore = $true
} catch {}

# =========================================================================
# Tabelle2 / Tabelle3 are otherwise unaffected by the diff (only namespace /
# ext metadata differences which aren't reachable through the object model).
# =========================================================================
